$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "complete MapData & subMap" -- the data table (rows 1..50 = MapData,
# row 51 begins the new sub-map entry) was missing its last vertex row.
# Duplicate the row immediately above (row 50) so the new row inherits
# the same cell formatting/number format as the rest of the table, then
# overwrite it with the new vertex data: id 50, name "RuanjianN",
# category 2, sub-category 1.
$ws.Rows(50).Copy() | Out-Null
$ws.Rows(51).Insert(-4121) | Out-Null
$ws.Rows(69).Delete() | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A51").Value = "50"
$ws.Range("B51").Value = "RuanjianN"
$ws.Range("C51").Value = "2"
$ws.Range("D51").Value = "1"

# The author's cursor ended up on D12 (no longer scrolled down to A20)
# when the workbook was saved.
$ws.Range("D12").Select() | Out-Null
